# Load_Profile_Range.xlsx update
# - Update the title string (date range) in A1
# - Update the load-profile demand values in column B (rows 3-26)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title text (date range changed from 2023-04-01..2023-04-10 to 2025-03-01..2025-03-03)
$ws.Range("A1").Value = "Total Load Profile for All Users for: 2025-03-01 to 2025-03-03 (Mthembanji)"

# Update demand values
$ws.Range("B3").Value = 0.4387077
$ws.Range("B4").Value = 0.4249936521
$ws.Range("B5").Value = 0.4318931247
$ws.Range("B6").Value = 0.42809348495
$ws.Range("B7").Value = 0.32199434665
$ws.Range("B8").Value = 0.3024936521
$ws.Range("B9").Value = 0.4139936521
$ws.Range("B10").Value = 0.5844936520999999
$ws.Range("B11").Value = 1.13354653675
$ws.Range("B12").Value = 1.32479653675
$ws.Range("B13").Value = 1.34379653675
$ws.Range("B14").Value = 1.42104653675
$ws.Range("B15").Value = 1.40779653675
$ws.Range("B16").Value = 1.50725
$ws.Range("B17").Value = 1.5722817274
$ws.Range("B18").Value = 1.4477182726
$ws.Range("B19").Value = 1.8672604107
$ws.Range("B20").Value = 2.0404983999
$ws.Range("B21").Value = 1.8135516456
$ws.Range("B22").Value = 1.3557847818
$ws.Range("B23").Value = 0.9550238095
$ws.Range("B24").Value = 0.8125238095
$ws.Range("B25").Value = 0.7765238095
$ws.Range("B26").Value = 0.7483097616000001
